$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 287705
$ws.Cells.Item(2, 4).Value = 367332540
$ws.Cells.Item(10, 3).Value = 109010
$ws.Cells.Item(10, 4).Value = 159864524
$ws.Cells.Item(12, 3).Value = 53840
$ws.Cells.Item(12, 4).Value = 77763562
$ws.Cells.Item(16, 3).Value = 3604
$ws.Cells.Item(16, 4).Value = 5121962
$ws.Cells.Item(20, 3).Value = 5393
$ws.Cells.Item(20, 4).Value = 7539526
$ws.Cells.Item(22, 3).Value = 70707
$ws.Cells.Item(22, 4).Value = 88595719
$ws.Cells.Item(28, 3).Value = 30445
$ws.Cells.Item(28, 4).Value = 44593492
$ws.Cells.Item(30, 3).Value = 10523
$ws.Cells.Item(30, 4).Value = 15159079
$ws.Cells.Item(35, 3).Value = 1511
$ws.Cells.Item(35, 4).Value = 2131079
$ws.Cells.Item(36, 3).Value = 89383
$ws.Cells.Item(36, 4).Value = 113069312
$ws.Cells.Item(42, 3).Value = 862
$ws.Cells.Item(42, 4).Value = 1270961
$ws.Cells.Item(44, 3).Value = 41743
$ws.Cells.Item(44, 4).Value = 61244082
$ws.Cells.Item(46, 3).Value = 8364
$ws.Cells.Item(46, 4).Value = 12012270
$ws.Cells.Item(48, 3).Value = 1276
$ws.Cells.Item(48, 4).Value = 1769466
$ws.Cells.Item(52, 3).Value = 62651
$ws.Cells.Item(52, 4).Value = 78796246
$ws.Cells.Item(58, 3).Value = 26276
$ws.Cells.Item(58, 4).Value = 38552791
$ws.Cells.Item(60, 3).Value = 5
$ws.Cells.Item(60, 4).Value = 7500
$ws.Cells.Item(61, 3).Value = 10128
$ws.Cells.Item(61, 4).Value = 14651985
$ws.Cells.Item(63, 3).Value = 1250
$ws.Cells.Item(63, 4).Value = 1744974
$ws.Cells.Item(68, 3).Value = 18206
$ws.Cells.Item(68, 4).Value = 23817317
$ws.Cells.Item(72, 3).Value = 6650
$ws.Cells.Item(72, 4).Value = 9729021
$ws.Cells.Item(74, 3).Value = 4526
$ws.Cells.Item(74, 4).Value = 6575350
$ws.Cells.Item(77, 3).Value = 128051
$ws.Cells.Item(77, 4).Value = 160087080
$ws.Cells.Item(83, 3).Value = 59427
$ws.Cells.Item(83, 4).Value = 87188231
$ws.Cells.Item(86, 3).Value = 27350
$ws.Cells.Item(86, 4).Value = 39579845
$ws.Cells.Item(88, 3).Value = 2471
$ws.Cells.Item(88, 4).Value = 3560581
$ws.Cells.Item(89, 3).Value = 2311
$ws.Cells.Item(89, 4).Value = 3259173
$ws.Cells.Item(90, 3).Value = 26328
$ws.Cells.Item(90, 4).Value = 35695222
$ws.Cells.Item(94, 3).Value = 6723
$ws.Cells.Item(94, 4).Value = 9907574
$ws.Cells.Item(96, 3).Value = 5906
$ws.Cells.Item(96, 4).Value = 8554912
$ws.Cells.Item(99, 3).Value = 390
$ws.Cells.Item(99, 4).Value = 563524
$ws.Cells.Item(100, 3).Value = 6305
$ws.Cells.Item(100, 4).Value = 8727434
$ws.Cells.Item(102, 3).Value = 1594
$ws.Cells.Item(102, 4).Value = 2343907
$ws.Cells.Item(104, 3).Value = 2118
$ws.Cells.Item(104, 4).Value = 3084017
$ws.Cells.Item(108, 3).Value = 128677
$ws.Cells.Item(108, 4).Value = 159311123
$ws.Cells.Item(114, 3).Value = 49477
$ws.Cells.Item(114, 4).Value = 72591216
$ws.Cells.Item(116, 3).Value = 24594
$ws.Cells.Item(116, 4).Value = 35634586
$ws.Cells.Item(120, 3).Value = 1877
$ws.Cells.Item(120, 4).Value = 2633978
$ws.Cells.Item(122, 3).Value = 415457
$ws.Cells.Item(122, 4).Value = 546044895
$ws.Cells.Item(129, 3).Value = 186246
$ws.Cells.Item(129, 4).Value = 273996871
$ws.Cells.Item(132, 3).Value = 156092
$ws.Cells.Item(132, 4).Value = 226893953
$ws.Cells.Item(135, 3).Value = 2393
$ws.Cells.Item(135, 4).Value = 3354320
$ws.Cells.Item(140, 3).Value = 38762
$ws.Cells.Item(140, 4).Value = 51904035
$ws.Cells.Item(146, 3).Value = 12821
$ws.Cells.Item(146, 4).Value = 18834013
$ws.Cells.Item(147, 3).Value = 3343
$ws.Cells.Item(147, 4).Value = 4827446
$ws.Cells.Item(153, 3).Value = 14932
$ws.Cells.Item(153, 4).Value = 19801300
$ws.Cells.Item(157, 3).Value = 6309
$ws.Cells.Item(157, 4).Value = 9193014
$ws.Cells.Item(159, 3).Value = 4213
$ws.Cells.Item(159, 4).Value = 6077661
$ws.Cells.Item(164, 3).Value = 10953
$ws.Cells.Item(164, 4).Value = 15831059
$ws.Cells.Item(165, 3).Value = 1403
$ws.Cells.Item(165, 4).Value = 2086078
$ws.Cells.Item(167, 3).Value = 29
$ws.Cells.Item(167, 4).Value = 43190
$ws.Cells.Item(169, 3).Value = 80013
$ws.Cells.Item(169, 4).Value = 100399221
$ws.Cells.Item(174, 3).Value = 617
$ws.Cells.Item(174, 4).Value = 910426
$ws.Cells.Item(176, 3).Value = 31780
$ws.Cells.Item(176, 4).Value = 46630495
$ws.Cells.Item(178, 3).Value = 11896
$ws.Cells.Item(178, 4).Value = 17198688
$ws.Cells.Item(180, 3).Value = 1133
$ws.Cells.Item(180, 4).Value = 1585120
$ws.Cells.Item(182, 3).Value = 1367
$ws.Cells.Item(182, 4).Value = 1919636
$ws.Cells.Item(184, 3).Value = 215938
$ws.Cells.Item(184, 4).Value = 268986047
$ws.Cells.Item(192, 3).Value = 81137
$ws.Cells.Item(192, 4).Value = 118985428
$ws.Cells.Item(195, 3).Value = 30233
$ws.Cells.Item(195, 4).Value = 43521672
$ws.Cells.Item(198, 3).Value = 4605
$ws.Cells.Item(198, 4).Value = 6563546
$ws.Cells.Item(201, 3).Value = 4012
$ws.Cells.Item(201, 4).Value = 5561992
$ws.Cells.Item(204, 3).Value = 236365
$ws.Cells.Item(204, 4).Value = 292975803
$ws.Cells.Item(213, 3).Value = 88238
$ws.Cells.Item(213, 4).Value = 129159169
$ws.Cells.Item(216, 3).Value = 46823
$ws.Cells.Item(216, 4).Value = 67721468
$ws.Cells.Item(219, 3).Value = 4230
$ws.Cells.Item(219, 4).Value = 5937630
$ws.Cells.Item(222, 3).Value = 4638
$ws.Cells.Item(222, 4).Value = 6403441
$ws.Cells.Item(225, 3).Value = 96538
$ws.Cells.Item(225, 4).Value = 121257723
$ws.Cells.Item(232, 3).Value = 46466
$ws.Cells.Item(232, 4).Value = 68108650
$ws.Cells.Item(234, 3).Value = 11248
$ws.Cells.Item(234, 4).Value = 16178493
$ws.Cells.Item(236, 3).Value = 1769
$ws.Cells.Item(236, 4).Value = 2537233
$ws.Cells.Item(238, 3).Value = 2131
$ws.Cells.Item(238, 4).Value = 2967463
$ws.Cells.Item(239, 3).Value = 231435
$ws.Cells.Item(239, 4).Value = 292523666
$ws.Cells.Item(241, 3).Value = 230
$ws.Cells.Item(241, 4).Value = 331342
$ws.Cells.Item(247, 3).Value = 88762
$ws.Cells.Item(247, 4).Value = 130149263
$ws.Cells.Item(250, 3).Value = 58894
$ws.Cells.Item(250, 4).Value = 85401122
$ws.Cells.Item(252, 3).Value = 2186
$ws.Cells.Item(252, 4).Value = 3084577
